# "Danh gia lan 3." -- add a third "tong quan" review bullet to the PA plan,
# and drop the now-obsolete "PA3" worksheet.

$wb = $excel.ActiveWorkbook

# --- 1. Remove the old "PA3" sheet -----------------------------------------
$wb.Worksheets("PA3").Delete()

# Make sure "PA" (now the first tab) is the active one, so the workbook
# doesn't keep pointing at a stale tab index after the deletion.
$wb.Worksheets("PA").Activate()

# --- 2. Append a "Danh gia tong quan lan 3" bullet to every PA row that
#        already ends with "...lan 2" (column E, rows 2-6) ------------------
$ws = $wb.Worksheets("PA")

$bulletChar = [string][char]0xB7   # "\u00b7" middle dot, same glyph as the
                                    # existing "." bullets typed with Symbol font
$newLine = "Đánh giá tổng quan lần 3"

$rows = @(2, 3, 4, 5, 6)
$newRowHeights = @{ 2 = 90; 3 = 75; 4 = 90; 5 = 90; 6 = 75 }

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 5)   # column E
    $oldText = $cell.Text
    $oldLen = $oldText.Length

    $fullText = $oldText + "`n" + $bulletChar + " " + $newLine
    $cell.Value2 = $fullText

    # Re-apply the bullet/body rich-text formatting across the *whole*
    # string (Value2 assignment flattens existing runs), then the new
    # trailing bullet gets the same Symbol/Arial split as all the others.
    $pos = 1
    $lines = $fullText -split "`n"
    for ($i = 0; $i -lt $lines.Count; $i++) {
        $line = $lines[$i]

        $bulletRun = $cell.Characters($pos, 1)
        $bulletRun.Font.Name = "Symbol"
        $bulletRun.Font.Size = 11
        $pos = $pos + 1

        $restLen = $line.Length - 1
        if ($i -lt ($lines.Count - 1)) {
            $restLen = $restLen + 1   # include the trailing "`n" in this run
        }
        $restRun = $cell.Characters($pos, $restLen)
        $restRun.Font.Name = "Arial"
        $restRun.Font.Size = 11
        $pos = $pos + $restLen
    }

    # The row now wraps one more line of text, so its cached height grows
    # by 15pt, matching Excel's own auto height for this wrapped style.
    $ws.Rows($r).RowHeight = $newRowHeights[$r]
}
